# All Country Files Saved And Formatted
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text tweak
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "MasterSheet RowNo."

# ---------------------------------------------------------------------------
# 2) Fill in the previously-missing TotalConfirmedNewCases (G) /
#    TotalNewDeaths (I) columns for the existing rows.
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = 1
$ws.Range("I2").Value = 0

$ws.Range("G3").Value = 2
$ws.Range("I3").Value = 0

$ws.Range("G4").Value = 0
$ws.Range("I4").Value = 0

$ws.Range("G5").Value = 0
$ws.Range("I5").Value = 0

$ws.Range("G6").Value = 0
$ws.Range("I6").Value = 0

$ws.Range("G7").Value = 0
$ws.Range("I7").Value = 0

$ws.Range("G8").Value = 2
$ws.Range("I8").Value = 1

$ws.Range("G9").Value = 0
$ws.Range("I9").Value = 0

$ws.Range("G10").Value = 0
$ws.Range("I10").Value = 0

# Correct the MasterSheet row number for row 10 (5215 -> 5214)
$ws.Range("L10").Value = 5214

# ---------------------------------------------------------------------------
# 3) Append the new data row (row 11)
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = 71
$ws.Range("B11").Value = 38
$ws.Range("C11").Value = "SUB-SAHARAN AFRICA                 "
$ws.Range("D11").Value = 43921
$ws.Range("E11").Value = "Cape Verde"
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "Imported cases only"
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 5417

# ---------------------------------------------------------------------------
# 4) Column widths: A..O all become 27 "characters" wide.
#    (ColumnWidth uses Excel's character-width units, which carry a fixed
#    +5/6 padding versus the raw OOXML <col width> value for the default
#    Calibri 11 font, hence the offset below.)
# ---------------------------------------------------------------------------
for ($c = 1; $c -le 15; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 26.166666666666668
}

# ---------------------------------------------------------------------------
# 5) Formatting: center (horizontal + vertical) every cell in A1:O11, and
#    additionally give column D a custom date format. Each combined style
#    is built once on a scratch cell and stamped onto the target cells via
#    copy / paste-special so only the exact styles we need get added to the
#    workbook's cellXfs table (no throw-away intermediate styles).
# ---------------------------------------------------------------------------
$tmplGeneral = $ws.Range("Z1")
$tmplGeneral.HorizontalAlignment = -4108
$tmplGeneral.VerticalAlignment = -4108
$tmplGeneral.Copy()
$ws.Range("A1:C11").PasteSpecial(-4122)
$ws.Range("E1:O11").PasteSpecial(-4122)
$tmplGeneral.Clear()

$tmplDate = $ws.Range("Z1")
$tmplDate.NumberFormat = "yyyy-mm-dd;"
$tmplDate.HorizontalAlignment = -4108
$tmplDate.VerticalAlignment = -4108
$tmplDate.Copy()
$ws.Range("D1:D11").PasteSpecial(-4122)
$tmplDate.Clear()

Write-Output "ok"
